# Actualización desde MV -datos-
# Appends 5 new daily rows (04-10-2021 .. 08-10-2021) to the "Diaria" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append: Date (text), "1 en 1", "3 en 2", "5 en 5"
$newRows = @(
    @("04-10-2021", 3.4,  3.16, 3.24),
    @("05-10-2021", 3.44, 3.19, 3.19),
    @("06-10-2021", 3.38, 3.19, 3.2),
    @("07-10-2021", 3.41, 3.2,  3.17),
    @("08-10-2021", 3.61, 3.27, 3.12)
)

# Find the last used row in column A and append right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

foreach ($row in $newRows) {
    $lastRow = $lastRow + 1

    # Force column A to remain plain text (so the date-like string
    # "dd-mm-yyyy" is not auto-converted into a date serial number), then
    # restore the "Normal" cell style so no extra formatting is left behind,
    # matching the plain (unstyled) text cells used by the rest of column A.
    $cellA = $ws.Cells.Item($lastRow, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($lastRow, 2).Value = $row[1]
    $ws.Cells.Item($lastRow, 3).Value = $row[2]
    $ws.Cells.Item($lastRow, 4).Value = $row[3]
}
